$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(70, 8).Value = 2023.75
$ws.Cells.Item(70, 9).Value = 1931.2222
$ws.Cells.Item(70, 11).Value = 5793.6666
$ws.Cells.Item(70, 13).Value = -5523.6666

$ws.Cells.Item(73, 8).Value = 2023.75
$ws.Cells.Item(73, 9).Value = 1931.2222
$ws.Cells.Item(73, 11).Value = 5793.6666
$ws.Cells.Item(73, 13).Value = -4857.6666

$ws.Cells.Item(88, 8).Value = 33782.668
$ws.Cells.Item(88, 9).Value = 350
$ws.Cells.Item(88, 10).Value = 50499
$ws.Cells.Item(88, 11).Value = 350
$ws.Cells.Item(88, 12).Value = 50499
$ws.Cells.Item(88, 13).Value = 56
$ws.Cells.Item(88, 14).Value = -51311

$ws.Cells.Item(91, 8).Value = 33782.668
$ws.Cells.Item(91, 9).Value = 350
$ws.Cells.Item(91, 10).Value = 50499
$ws.Cells.Item(91, 11).Value = 350
$ws.Cells.Item(91, 12).Value = 50499
$ws.Cells.Item(91, 13).Value = 1054
$ws.Cells.Item(91, 14).Value = -53307

$ws.Cells.Item(93, 8).Value = 15015.333
$ws.Cells.Item(93, 10).Value = 15015.333
$ws.Cells.Item(93, 12).Value = 15015.333
$ws.Cells.Item(93, 14).Value = -20007.333

$ws.Cells.Item(100, 8).Value = 4256
$ws.Cells.Item(100, 10).Value = 3100
$ws.Cells.Item(100, 12).Value = 3100
$ws.Cells.Item(100, 14).Value = -4182

$ws.Cells.Item(116, 8).Value = 4000
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).ClearContents()

$ws.Cells.Item(138, 8).Value = 2979.9412
$ws.Cells.Item(138, 9).Value = 1406.5555
$ws.Cells.Item(138, 10).Value = 4750
$ws.Cells.Item(138, 11).Value = 4219.666499999999
$ws.Cells.Item(138, 12).Value = 14250
$ws.Cells.Item(138, 13).Value = 920.3335000000006
$ws.Cells.Item(138, 14).Value = -24530

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 3929.9348
$ws.Cells.Item(32, 9).Value = 1219.475
$ws.Cells.Item(32, 11).Value = 1219.475
$ws.Cells.Item(32, 13).Value = -932.4749999999999

$ws.Cells.Item(45, 8).Value = 1739.8
$ws.Cells.Item(45, 9).Value = 1674.75
$ws.Cells.Item(45, 11).Value = 1674.75
$ws.Cells.Item(45, 13).Value = -1297.75

$ws.Cells.Item(97, 8).Value = 602.75
$ws.Cells.Item(97, 9).Value = 592.3333
$ws.Cells.Item(97, 11).Value = 592.3333
$ws.Cells.Item(97, 13).Value = -96.33330000000001

$ws.Cells.Item(122, 8).Value = 2873.25
$ws.Cells.Item(122, 9).Value = 2914.3333
$ws.Cells.Item(122, 10).Value = 2750
$ws.Cells.Item(122, 11).Value = 8742.999899999999
$ws.Cells.Item(122, 12).Value = 8250
$ws.Cells.Item(122, 13).Value = -6292.999899999999
$ws.Cells.Item(122, 14).Value = -13150

$ws.Cells.Item(132, 8).Value = 4386.2856
$ws.Cells.Item(132, 9).Value = 4153.647
$ws.Cells.Item(132, 10).Value = 5375
$ws.Cells.Item(132, 11).Value = 12460.941
$ws.Cells.Item(132, 12).Value = 16125
$ws.Cells.Item(132, 13).Value = -9930.940999999999
$ws.Cells.Item(132, 14).Value = -21185

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 4129.143
$ws.Cells.Item(86, 9).Value = 3984
$ws.Cells.Item(86, 11).Value = 3984
$ws.Cells.Item(86, 13).Value = -2861

$ws.Cells.Item(89, 8).Value = 4129.143
$ws.Cells.Item(89, 9).Value = 3984
$ws.Cells.Item(89, 11).Value = 19920
$ws.Cells.Item(89, 13).Value = -14304

$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).ClearContents()
$ws.Cells.Item(92, 14).Value = 0

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 306.41666
$ws.Cells.Item(22, 9).Value = 153
$ws.Cells.Item(22, 10).Value = 613.25
$ws.Cells.Item(22, 11).Value = 153
$ws.Cells.Item(22, 12).Value = 613.25
$ws.Cells.Item(22, 13).Value = 197
$ws.Cells.Item(22, 14).Value = -1313.25

$ws.Cells.Item(58, 8).Value = 1756.6522
$ws.Cells.Item(58, 10).Value = 1749.875
$ws.Cells.Item(58, 12).Value = 1749.875
$ws.Cells.Item(58, 14).Value = -2155.875

$ws.Cells.Item(68, 8).Value = 29997.727
$ws.Cells.Item(68, 10).Value = 29997.727
$ws.Cells.Item(68, 12).Value = 29997.727
$ws.Cells.Item(68, 14).Value = -31495.727

$ws.Cells.Item(71, 8).Value = 29997.727
$ws.Cells.Item(71, 10).Value = 29997.727
$ws.Cells.Item(71, 12).Value = 89993.181
$ws.Cells.Item(71, 14).Value = -97481.181

$ws.Cells.Item(74, 8).Value = 29089.092
$ws.Cells.Item(74, 10).Value = 29089.092
$ws.Cells.Item(74, 12).Value = 29089.092
$ws.Cells.Item(74, 14).Value = -30837.092

$ws.Cells.Item(77, 8).Value = 29089.092
$ws.Cells.Item(77, 10).Value = 29089.092
$ws.Cells.Item(77, 12).Value = 87267.276
$ws.Cells.Item(77, 14).Value = -96003.276

$ws.Cells.Item(93, 8).Value = 10703.5
$ws.Cells.Item(93, 9).Value = 10703.5
$ws.Cells.Item(93, 11).Value = 10703.5
$ws.Cells.Item(93, 13).Value = -8831.5

$ws.Cells.Item(122, 8).Value = 2077
$ws.Cells.Item(122, 9).Value = 1781.75
$ws.Cells.Item(122, 11).Value = 5345.25
$ws.Cells.Item(122, 13).Value = -2895.25

$ws.Cells.Item(132, 8).Value = 3808
$ws.Cells.Item(132, 9).Value = 3637
$ws.Cells.Item(132, 10).Value = 4150
$ws.Cells.Item(132, 11).Value = 10911
$ws.Cells.Item(132, 12).Value = 12450
$ws.Cells.Item(132, 13).Value = -8381
$ws.Cells.Item(132, 14).Value = -17510

$ws.Cells.Item(136, 8).Value = 1756.6522
$ws.Cells.Item(136, 10).Value = 1749.875
$ws.Cells.Item(136, 12).Value = 5249.625
$ws.Cells.Item(136, 14).Value = -10349.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 1542.52
$ws.Cells.Item(4, 9).Value = 1387.5264
$ws.Cells.Item(4, 11).Value = 4162.5792
$ws.Cells.Item(4, 13).Value = -4050.5792

$ws.Cells.Item(38, 8).Value = 207.28572
$ws.Cells.Item(38, 10).Value = 95.2
$ws.Cells.Item(38, 12).Value = 285.6
$ws.Cells.Item(38, 14).Value = -979.6

$ws.Cells.Item(59, 8).Value = 0
$ws.Cells.Item(59, 10).Value = 0
$ws.Cells.Item(59, 12).ClearContents()
$ws.Cells.Item(59, 14).Value = 0

$ws.Cells.Item(97, 8).Value = 1229.3334
$ws.Cells.Item(97, 10).Value = 1229.3334
$ws.Cells.Item(97, 12).Value = 3688.0002
$ws.Cells.Item(97, 14).Value = -4680.0002

$ws.Cells.Item(131, 8).Value = 591.3333
$ws.Cells.Item(131, 9).Value = 591.3333
$ws.Cells.Item(131, 10).Value = 0
$ws.Cells.Item(131, 11).Value = 1773.9999
$ws.Cells.Item(131, 12).Value = 0
$ws.Cells.Item(131, 13).ClearContents()
$ws.Cells.Item(131, 14).Value = 3266.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 5000
$ws.Cells.Item(5, 9).Value = 5000
$ws.Cells.Item(5, 11).Value = 5000
$ws.Cells.Item(5, 13).Value = -4888

$ws.Cells.Item(13, 8).Value = 318.5
$ws.Cells.Item(13, 9).Value = 150
$ws.Cells.Item(13, 10).Value = 402.75
$ws.Cells.Item(13, 11).Value = 150
$ws.Cells.Item(13, 12).Value = 402.75
$ws.Cells.Item(13, 13).Value = -11
$ws.Cells.Item(13, 14).Value = -680.75

$ws.Cells.Item(80, 8).Value = 1846.6
$ws.Cells.Item(80, 9).Value = 1846.6
$ws.Cells.Item(80, 11).Value = 1846.6
$ws.Cells.Item(80, 13).Value = -848.5999999999999

$ws.Cells.Item(83, 8).Value = 1846.6
$ws.Cells.Item(83, 9).Value = 1846.6
$ws.Cells.Item(83, 11).Value = 9233
$ws.Cells.Item(83, 13).Value = -4241

$ws.Cells.Item(96, 8).Value = 29475
$ws.Cells.Item(96, 10).Value = 29475
$ws.Cells.Item(96, 12).Value = 29475
$ws.Cells.Item(96, 14).Value = -34967

$ws.Cells.Item(102, 8).Value = 12131.75
$ws.Cells.Item(102, 9).Value = 3224.182
$ws.Cells.Item(102, 11).Value = 3224.182
$ws.Cells.Item(102, 13).Value = -1602.182

$ws.Cells.Item(107, 8).Value = 500.625
$ws.Cells.Item(107, 9).Value = 486.2
$ws.Cells.Item(107, 10).Value = 524.6667
$ws.Cells.Item(107, 11).Value = 486.2
$ws.Cells.Item(107, 12).Value = 524.6667
$ws.Cells.Item(107, 13).Value = 1433.8
$ws.Cells.Item(107, 14).Value = -4364.6667

$ws.Cells.Item(122, 8).Value = 3385.125
$ws.Cells.Item(122, 9).Value = 3100
$ws.Cells.Item(122, 11).Value = 9300
$ws.Cells.Item(122, 13).Value = -6850

$ws.Cells.Item(126, 8).Value = 9398
$ws.Cells.Item(126, 9).Value = 6097
$ws.Cells.Item(126, 11).Value = 18291
$ws.Cells.Item(126, 13).Value = -15821

$ws.Cells.Item(131, 8).Value = 99000
$ws.Cells.Item(131, 10).Value = 99000
$ws.Cells.Item(131, 12).Value = 99000
$ws.Cells.Item(131, 14).Value = -109080

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2023.9166
$ws.Cells.Item(16, 9).Value = 2023.9166
$ws.Cells.Item(16, 11).Value = 2023.9166
$ws.Cells.Item(16, 13).Value = -1853.9166

$ws.Cells.Item(22, 8).Value = 1333.3334
$ws.Cells.Item(22, 9).Value = 1000
$ws.Cells.Item(22, 10).Value = 2000
$ws.Cells.Item(22, 11).Value = 1000
$ws.Cells.Item(22, 12).Value = 2000
$ws.Cells.Item(22, 13).Value = -705
$ws.Cells.Item(22, 14).Value = -2590

$ws.Cells.Item(27, 8).Value = 1333.3334
$ws.Cells.Item(27, 9).Value = 1000
$ws.Cells.Item(27, 10).Value = 2000
$ws.Cells.Item(27, 11).Value = 1000
$ws.Cells.Item(27, 12).Value = 2000
$ws.Cells.Item(27, 13).Value = -893
$ws.Cells.Item(27, 14).Value = -2214

$ws.Cells.Item(46, 8).Value = 4981.6
$ws.Cells.Item(46, 9).Value = 5591.5713
$ws.Cells.Item(46, 10).Value = 3558.3333
$ws.Cells.Item(46, 11).Value = 5591.5713
$ws.Cells.Item(46, 12).Value = 3558.3333
$ws.Cells.Item(46, 13).Value = -5403.5713
$ws.Cells.Item(46, 14).Value = -3934.3333

$ws.Cells.Item(55, 8).Value = 320.72726
$ws.Cells.Item(55, 10).Value = 354.66666
$ws.Cells.Item(55, 12).Value = 354.66666
$ws.Cells.Item(55, 14).Value = -700.66666

$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 9).Value = 0
$ws.Cells.Item(93, 11).Value = 0
$ws.Cells.Item(93, 13).ClearContents()

$ws.Cells.Item(122, 8).Value = 3140.6667
$ws.Cells.Item(122, 9).Value = 3287.2856
$ws.Cells.Item(122, 10).Value = 2627.5
$ws.Cells.Item(122, 11).Value = 9861.856800000001
$ws.Cells.Item(122, 12).Value = 7882.5
$ws.Cells.Item(122, 13).Value = -7411.856800000001
$ws.Cells.Item(122, 14).Value = -12782.5

$ws.Cells.Item(132, 8).Value = 2800
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 13).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 83799.8
$ws.Cells.Item(2, 10).Value = 102749.5
$ws.Cells.Item(2, 12).Value = 102749.5
$ws.Cells.Item(2, 14).Value = -102973.5

$ws.Cells.Item(107, 8).Value = 1400
$ws.Cells.Item(107, 9).Value = 966.8333
$ws.Cells.Item(107, 10).Value = 2049.75
$ws.Cells.Item(107, 11).Value = 2900.4999
$ws.Cells.Item(107, 12).Value = 6149.25
$ws.Cells.Item(107, 13).Value = -980.4998999999998
$ws.Cells.Item(107, 14).Value = -9989.25

$ws.Cells.Item(132, 8).Value = 0
$ws.Cells.Item(132, 10).Value = 0
$ws.Cells.Item(132, 12).ClearContents()
$ws.Cells.Item(132, 14).Value = 0
